# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / handoff / handback timestamp
# cells to reflect a fresh report generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date for b135cdc8-... (row 4)
$wsOverview.Range("G4").Value = "2016-08-25 22:44:40"

# zh-cn sheet: Correspond Handoff / Handback DateTime for b135cdc8-... (row 4)
$wsZhCn.Range("H4").Value = "2016-08-25 22:44:35"
$wsZhCn.Range("K4").Value = "2016-08-25 22:45:06"

# de-de sheet: Correspond Handback DateTime for b135cdc8-... (row 4)
$wsDeDe.Range("K4").Value = "2016-08-25 22:45:20"
